{"js": "// Remove the \"Expand Timeframe\" bullet point paragraph from the\n// \"Future Enhancements\" list (keep \"Advanced Analytics\" and the heading).\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst target = \"Expand Timeframe:\";\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  if (para.text && para.text.indexOf(target) !== -1) {\n    para.delete();\n  }\n}\nawait context.sync();\n", "ps1": "# Remove the \"Expand Timeframe\" bullet point paragraph from the\n# \"Future Enhancements\" list (keep \"Advanced Analytics\" and the heading).\n$d = $word.ActiveDocument\n\n$searchText = \"Expand Timeframe:\"\n\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$rng.Find.Forward = $true\n$rng.Find.Wrap = 1\n\nwhile ($rng.Find.Execute($searchText)) {\n  # Expand the found hit to its whole paragraph (including the\n  # paragraph mark) so the entire bullet point is removed.\n  [void]$rng.Expand(4)\n  $rng.Delete()\n\n  # Re-scope the search range to the remainder of the document so we\n  # don't loop forever / re-find content that no longer exists.\n  $rng = $d.Content\n  $rng.Find.ClearFormatting()\n  $rng.Find.Forward = $true\n  $rng.Find.Wrap = 1\n}\n"}
